$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price updates (column D). Leading apostrophe forces text storage so the
# values stay strings (e.g. "0.05620") instead of being coerced to numbers.
$ws.Cells.Item(2, 4).Value = "'247.24"
$ws.Cells.Item(3, 4).Value = "'22.41"
$ws.Cells.Item(4, 4).Value = "'5.485"
$ws.Cells.Item(5, 4).Value = "'0.05620"
$ws.Cells.Item(6, 4).Value = "'6.456"
$ws.Cells.Item(8, 4).Value = "'1.044"
$ws.Cells.Item(9, 4).Value = "'0.1420"
$ws.Cells.Item(10, 4).Value = "'0.07314"
$ws.Cells.Item(11, 4).Value = "'0.03200"
$ws.Cells.Item(12, 4).Value = "'0.02921"
$ws.Cells.Item(13, 4).Value = "'0.09249"
$ws.Cells.Item(14, 4).Value = "'0.001663"
$ws.Cells.Item(15, 4).Value = "'3.208"
$ws.Cells.Item(16, 4).Value = "'0.04729"

# Row 17-24: "One" moved up from row 24 to row 17, shifting TigerCash,
# HotbitToken, BitKan, NitroEx, LEO, GateToken and BTSEToken down by one
# row (with refreshed prices/volume labels for each).
$ws.Cells.Item(17, 2).Value = "One"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Cells.Item(17, 4).Value = "'0.0005821"
$ws.Cells.Item(17, 5).Value = "16OneONE"
$ws.Cells.Item(18, 2).Value = "TigerCash"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(18, 4).Value = "'0.006438"
$ws.Cells.Item(18, 5).Value = "17TigerCashTCH"
$ws.Cells.Item(19, 2).Value = "HotbitToken"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Cells.Item(19, 4).Value = "'0.005065"
$ws.Cells.Item(19, 5).Value = "18HotbitTokenHTB"
$ws.Cells.Item(20, 2).Value = "BitKan"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Cells.Item(20, 4).Value = "'0.001057"
$ws.Cells.Item(20, 5).Value = "19BitKanKAN"
$ws.Cells.Item(21, 2).Value = "NitroEx"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Cells.Item(21, 4).Value = "'0.0001502"
$ws.Cells.Item(21, 5).Value = "20NitroExNTX"
$ws.Cells.Item(22, 2).Value = "LEO"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(22, 4).Value = "'3.981"
$ws.Cells.Item(22, 5).Value = "21LEOLEO"
$ws.Cells.Item(23, 2).Value = "GateToken"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(23, 4).Value = "'3.380"
$ws.Cells.Item(23, 5).Value = "22GateTokenGT"
$ws.Cells.Item(24, 2).Value = "BTSEToken"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Cells.Item(24, 4).Value = "'2.128"
$ws.Cells.Item(24, 5).Value = "23BTSETokenBTSE"
$ws.Cells.Item(25, 4).Value = "'0.3320"
$ws.Cells.Item(26, 5).Value = "25ProBitTokenPROBBestin24h"

$ws.Cells.Item(40, 4).Value = "'0.04156"

# Row 41 and 43 swap (KickToken <-> BKEXToken); row 42 (CEJI) stays put
# but its price is refreshed.
$ws.Cells.Item(41, 2).Value = "BKEXToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Cells.Item(41, 4).Value = "'0.1041"
$ws.Cells.Item(41, 5).Value = "40BKEXTokenBKK"
$ws.Cells.Item(42, 4).Value = "'0.002974"
$ws.Cells.Item(43, 2).Value = "KickToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Cells.Item(43, 4).Value = "'0.003247"
$ws.Cells.Item(43, 5).Value = "42KickTokenKICK"
$ws.Cells.Item(44, 4).Value = "'0.009113"
$ws.Cells.Item(45, 4).Value = "'0.00005657"
$ws.Cells.Item(48, 4).Value = "'0.01558"
